# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update price column (D) values for the two price tables
$ws.Range("D14").Value = 248.739
$ws.Range("D15").Value = 381.532
$ws.Range("D38").Value = 457.837
$ws.Range("D39").Value = 491.531
